# "add export and formate" -- expand the single-row header sheet from 11
# columns (A:K) out to 26 columns (A:Z), inserting several new header
# columns, widening every column to fit its new header text, bumping the
# header row height, and moving the selection/scroll position over to the
# newly added columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header row (A1:Z1), replacing/augmenting the old 11 headers.
# ---------------------------------------------------------------------
$headers = @(
    "Fleet Code",
    "Vehicle Number",
    "Agent Name",
    "Insurance Company",
    "Insurance Type",
    "Insured's Name",
    "Insurance Policy Number",
    "Insurance Amount ",
    "Insurance Total Amount",
    "Insurance Prev Policy No",
    "NCB Discount",
    "Hypothecation Agreement",
    "Payment Mode",
    "Pay Number",
    "Pay Date",
    "Pay Bank",
    "Pay Branch",
    "Valid From",
    "Valid Till",
    "Engine No.",
    "Chassis No",
    "Manufacture Year",
    "Type Of Body",
    "Type Of Fuel",
    "Seating Capacity(including Driver)",
    "Cubic Capacity"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Header row height.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15

# ---------------------------------------------------------------------
# 3. Per-column widths (character units), columns A-Z.
# ---------------------------------------------------------------------
$colWidths = @(
    25.427043269230772,
    21.141586538461542,
    27.99831730769231,
    40.140444711538464,
    20.855889423076924,
    26.85552884615385,
    30.56959134615385,
    30.283894230769235,
    30.56959134615385,
    31.283834134615386,
    20.427343750000002,
    31.71237980769231,
    16.998978365384616,
    21.71298076923077,
    17.99891826923077,
    20.284495192307695,
    26.28413461538462,
    17.42752403846154,
    16.141887019230772,
    15.856189903846156,
    22.284375,
    29.28395432692308,
    17.141826923076923,
    26.14128605769231,
    38.426262019230776,
    19.855949519230773
)

for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# ---------------------------------------------------------------------
# 4. View state: scroll right so column J is the leftmost visible
#    column, and select N2 (one row below the new "Pay Number" header).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("N2").Select() | Out-Null
